$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.740.19'
$ws.Range("D3").Value = '1.657.27'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.13'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3832'
$ws.Range("E7").Value = '  +0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3611'
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.16'
$ws.Range("E9").Value = '  -1.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08197'
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.230'
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.47'
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.452'
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.440'
$ws.Range("E15").Value = '  +1.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001225'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '1.656.56'
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.44'
$ws.Range("E18").Value = '  +2.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07024'
$ws.Range("E19").Value = '  +0.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.800'
$ws.Range("E20").Value = '  +3.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.59'
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("E23").Value = '  +1.73%  '
$ws.Range("D24").Value = '23.733.31'
$ws.Range("E24").Value = '  +1.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.502'
$ws.Range("E25").Value = '  -1.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.024'
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.25'
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '154.47'
$ws.Range("E28").Value = '  +1.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.241'
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.05'
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("B31").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C31").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D31").Value = '1.840.21'
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.219'
$ws.Range("E32").Value = '  +9.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.251'
$ws.Range("E33").Value = '  +4.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.99'
$ws.Range("E34").Value = '  +4.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.055'
$ws.Range("E35").Value = '  -3.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02806'
$ws.Range("E36").Value = '  +1.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2507'
$ws.Range("E37").Value = '  -0.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08801'
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("E39").Value = '  +1.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06994'
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.99'
$ws.Range("E41").Value = '  +5.62%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6995'
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("E43").Value = '  -1.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.97'
$ws.Range("E44").Value = '  +2.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6513'
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.302'
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.963'
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07908'
$ws.Range("E49").Value = '  -0.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.14'
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.179'
$ws.Range("E51").Value = '  -1.26%  '
